{"js": "// Replace the division problems in the table with their new values, as\n// described by the diff. Each old string is unique in the document, so a\n// simple search-and-replace per pair is safe and unambiguous.\nconst replacements = [\n  [\"28\u00f77=\", \"84\u00f75=\"],\n  [\"34\u00f79=\", \"90\u00f76=\"],\n  [\"96\u00f79=\", \"76\u00f76=\"],\n  [\"59\u00f76=\", \"28\u00f79=\"],\n  [\"96\u00f78=\", \"65\u00f74=\"],\n  [\"42\u00f72=\", \"29\u00f75=\"],\n  [\"11\u00f75=\", \"77\u00f77=\"],\n  [\"12\u00f77=\", \"27\u00f79=\"],\n  [\"75\u00f79=\", \"45\u00f79=\"],\n  [\"43\u00f72=\", \"14\u00f73=\"],\n  [\"30\u00f72=\", \"77\u00f76=\"],\n  [\"44\u00f76=\", \"17\u00f73=\"],\n  [\"45\u00f76=\", \"93\u00f79=\"],\n  [\"25\u00f76=\", \"89\u00f72=\"],\n  [\"97\u00f76=\", \"18\u00f79=\"],\n  [\"90\u00f73=\", \"57\u00f79=\"],\n  [\"85\u00f78=\", \"58\u00f76=\"],\n  [\"94\u00f76=\", \"81\u00f79=\"],\n  [\"49\u00f74=\", \"51\u00f78=\"],\n  [\"12\u00f79=\", \"38\u00f73=\"],\n  [\"62\u00f77=\", \"39\u00f77=\"],\n  [\"88\u00f76=\", \"49\u00f72=\"],\n  [\"81\u00f74=\", \"96\u00f74=\"],\n  [\"34\u00f77=\", \"76\u00f74=\"],\n  [\"50\u00f76=\", \"80\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the division problems in the table with their new values, as\n# described by the diff. Each old string is unique in the document, so a\n# simple Find/Replace per pair is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"28\u00f77=\", \"84\u00f75=\"),\n  @(\"34\u00f79=\", \"90\u00f76=\"),\n  @(\"96\u00f79=\", \"76\u00f76=\"),\n  @(\"59\u00f76=\", \"28\u00f79=\"),\n  @(\"96\u00f78=\", \"65\u00f74=\"),\n  @(\"42\u00f72=\", \"29\u00f75=\"),\n  @(\"11\u00f75=\", \"77\u00f77=\"),\n  @(\"12\u00f77=\", \"27\u00f79=\"),\n  @(\"75\u00f79=\", \"45\u00f79=\"),\n  @(\"43\u00f72=\", \"14\u00f73=\"),\n  @(\"30\u00f72=\", \"77\u00f76=\"),\n  @(\"44\u00f76=\", \"17\u00f73=\"),\n  @(\"45\u00f76=\", \"93\u00f79=\"),\n  @(\"25\u00f76=\", \"89\u00f72=\"),\n  @(\"97\u00f76=\", \"18\u00f79=\"),\n  @(\"90\u00f73=\", \"57\u00f79=\"),\n  @(\"85\u00f78=\", \"58\u00f76=\"),\n  @(\"94\u00f76=\", \"81\u00f79=\"),\n  @(\"49\u00f74=\", \"51\u00f78=\"),\n  @(\"12\u00f79=\", \"38\u00f73=\"),\n  @(\"62\u00f77=\", \"39\u00f77=\"),\n  @(\"88\u00f76=\", \"49\u00f72=\"),\n  @(\"81\u00f74=\", \"96\u00f74=\"),\n  @(\"34\u00f77=\", \"76\u00f74=\"),\n  @(\"50\u00f76=\", \"80\u00f76=\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
